$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'Datos actualizados a 4 de Abril de 2020 a las 01:50'
$ws.Range("B4").Value = 276931
$ws.Range("C4").Value = 32054
$ws.Range("E4").Value = 257272
$ws.Range("G4").Value = 1320
$ws.Range("H4").Value = 7391
$ws.Range("A42").Value = 'Panama'
$ws.Range("B42").Value = 1673
$ws.Range("C42").Value = 198
$ws.Range("D42").Value = 10
$ws.Range("E42").Value = 1622
$ws.Range("F42").Value = 50
$ws.Range("G42").Value = 4
$ws.Range("H42").Value = 41
$ws.Range("A43").Value = 'Finlandia'
$ws.Range("B43").Value = 1615
$ws.Range("C43").Value = 97
$ws.Range("D43").Value = 300
$ws.Range("E43").Value = 1295
$ws.Range("F43").Value = 72
$ws.Range("G43").Value = 1
$ws.Range("H43").Value = 20
$ws.Range("A44").Value = 'Grecia'
$ws.Range("B44").Value = 1613
$ws.Range("C44").Value = 69
$ws.Range("D44").Value = 78
$ws.Range("E44").Value = 1472
$ws.Range("F44").Value = 92
$ws.Range("G44").Value = 10
$ws.Range("H44").Value = 63
$ws.Range("A45").Value = 'Peru'
$ws.Range("B45").Value = 1595
$ws.Range("C45").Value = 181
$ws.Range("D45").Value = 537
$ws.Range("E45").Value = 997
$ws.Range("F45").Value = 51
$ws.Range("G45").Value = 6
$ws.Range("H45").Value = 61
$ws.Range("A46").Value = 'Mexico'
$ws.Range("B46").Value = 1510
$ws.Range("C46").Value = 132
$ws.Range("D46").Value = 633
$ws.Range("E46").Value = 827
$ws.Range("F46").Value = 1
$ws.Range("G46").Value = 13
$ws.Range("H46").Value = 50
$ws.Range("A47").Value = 'Sudafrica'
$ws.Range("B47").Value = 1505
$ws.Range("C47").Value = 43
$ws.Range("D47").Value = 95
$ws.Range("E47").Value = 1401
$ws.Range("F47").Value = 7
$ws.Range("G47").Value = 4
$ws.Range("H47").Value = 9
$ws.Range("A48").Value = 'Republica Dominicana'
$ws.Range("B48").Value = 1488
$ws.Range("C48").Value = 108
$ws.Range("D48").Value = 16
$ws.Range("E48").Value = 1404
$ws.Range("F48").Value = 147
$ws.Range("H48").Value = 68
$ws.Range("A49").Value = 'Serbia'
$ws.Range("B49").Value = 1476
$ws.Range("C49").Value = 305
$ws.Range("D49").Value = 54
$ws.Range("E49").Value = 1383
$ws.Range("F49").Value = 81
$ws.Range("G49").Value = 8
$ws.Range("H49").Value = 39
$ws.Range("A51").Value = 'Argentina'
$ws.Range("B51").Value = 1353
$ws.Range("C51").Value = 88
$ws.Range("D51").Value = 266
$ws.Range("E51").Value = 1045
$ws.Range("F51").Value = 0
$ws.Range("H51").Value = 42
$ws.Range("A52").Value = 'Colombia'
$ws.Range("B52").Value = 1267
$ws.Range("C52").Value = 106
$ws.Range("D52").Value = 55
$ws.Range("E52").Value = 1187
$ws.Range("F52").Value = 50
$ws.Range("G52").Value = 6
$ws.Range("H52").Value = 25
$ws.Range("F75").Value = 26
$ws.Range("B86").Value = 386
$ws.Range("C86").Value = 36
$ws.Range("D86").Value = 86
$ws.Range("E86").Value = 296
$ws.Range("B97").Value = 239
$ws.Range("C97").Value = 6
$ws.Range("E97").Value = 154
$ws.Range("B149").Value = 29
$ws.Range("C149").Value = 1
$ws.Range("D149").Value = 1
$ws.Range("E152").Value = 13
$ws.Range("G152").Value = 3
$ws.Range("H152").Value = 4
$ws.Range("A169").Value = 'Dominica'
$ws.Range("C169").Value = 2
$ws.Range("D169").Value = 0
$ws.Range("E169").Value = 14
$ws.Range("A170").Value = 'Mongolia'
$ws.Range("D170").Value = 2
$ws.Range("E170").Value = 12
$ws.Range("A171").Value = 'Namibia'
$ws.Range("B171").Value = 14
$ws.Range("D171").Value = 3
$ws.Range("E171").Value = 11
$ws.Range("A172").Value = 'Santa Lucia'
$ws.Range("B172").Value = 13
$ws.Range("C172").Value = 0
$ws.Range("D172").Value = 1
$ws.Range("F172").Value = 0
$ws.Range("A173").Value = 'Granada'
$ws.Range("C173").Value = 2
$ws.Range("F173").Value = 2
$ws.Range("A175").Value = 'Laos'
$ws.Range("A176").Value = 'Mozambique'
$ws.Range("A188").Value = 'Santa Sede'
$ws.Range("A190").Value = 'Fiyi'
$ws.Range("A191").Value = 'San Vicente y las Granadinas'
$ws.Range("C191").Value = 5
$ws.Range("A192").Value = 'Somalia'
$ws.Range("B192").Value = 7
$ws.Range("C192").Value = 2
$ws.Range("D192").Value = 1
$ws.Range("E192").Value = 6
$ws.Range("H192").Value = 0
$ws.Range("A193").Value = 'Cabo Verde'
$ws.Range("D193").Value = 0
$ws.Range("H193").Value = 1
$ws.Range("A194").Value = 'Nepal'
$ws.Range("A195").Value = 'San Bartolome'
$ws.Range("D195").Value = 1
$ws.Range("E195").Value = 5
$ws.Range("H195").Value = 0
$ws.Range("A196").Value = 'Mauritania'
$ws.Range("B196").Value = 6
$ws.Range("D196").Value = 2
$ws.Range("E196").Value = 3
$ws.Range("H196").Value = 1
$ws.Range("A197").Value = 'Islas Turcas y Caicos'
$ws.Range("E197").Value = 5
$ws.Range("H197").Value = 0
$ws.Range("A198").Value = 'Nicaragua'
$ws.Range("D198").Value = 0
$ws.Range("E198").Value = 4
$ws.Range("H198").Value = 1
$ws.Range("A199").Value = 'Butan'
$ws.Range("B199").Value = 5
$ws.Range("C199").Value = 0
$ws.Range("D199").Value = 2
$ws.Range("E199").Value = 3
$ws.Range("A200").Value = 'Belice'
$ws.Range("C200").Value = 1
$ws.Range("E200").Value = 4
$ws.Range("H200").Value = 0
$ws.Range("A201").Value = 'Botsuana'
$ws.Range("D201").Value = 0
$ws.Range("E201").Value = 3
$ws.Range("A202").Value = 'Gambia'
$ws.Range("B202").Value = 4
$ws.Range("D202").Value = 2
$ws.Range("E202").Value = 1
$ws.Range("H202").Value = 1
$ws.Range("A205").Value = 'Islas Virgenes Britanicas'
$ws.Range("A206").Value = 'Burundi'
$ws.Range("C206").Value = 0
$ws.Range("D206").Value = 0
$ws.Range("E206").Value = 3
$ws.Range("A207").Value = 'Bonaire, San Eustaquio y Saba'
$ws.Range("A208").Value = 'Sierra Leona'
